$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '57.980.76'
Set-TextValue 'E2' '  -1.71%  '
Set-TextValue 'D3' '2.451.14'
Set-TextValue 'E3' '  -3.58%  '
Set-TextValue 'E4' '  +0.18%  '
Set-TextValue 'D5' '524.61'
Set-TextValue 'E5' '  -0.32%  '
Set-TextValue 'E6' '  -2.98%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  -0.21%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'D9' '0.0978'
Set-TextValue 'E9' '  -1.08%  '
Set-TextValue 'E10' '  -1.93%  '
Set-TextValue 'E11' '  -4.00%  '
Set-TextValue 'E12' '  -3.80%  '
Set-TextValue 'D13' '2.886.45'
Set-TextValue 'E13' '  -3.50%  '
Set-TextValue 'D14' '57.910.13'
Set-TextValue 'E14' '  -1.71%  '
Set-TextValue 'E15' '  -3.68%  '
Set-TextValue 'E16' '  -2.56%  '
Set-TextValue 'D17' '2.452.95'
Set-TextValue 'E17' '  -3.50%  '
Set-TextValue 'D18' '10.39'
Set-TextValue 'E18' '  -3.17%  '
Set-TextValue 'E19' '  -1.86%  '
Set-TextValue 'D20' '312.05'
Set-TextValue 'E20' '  -3.61%  '
Set-TextValue 'D21' '6.15'
Set-TextValue 'E21' '  +0.27%  '
Set-TextValue 'D22' '0.999'
Set-TextValue 'E22' '  -0.09%  '
Set-TextValue 'D23' '64.86'
Set-TextValue 'E23' '  -0.47%  '
Set-TextValue 'D24' '0.402'
Set-TextValue 'E24' '  -2.07%  '
Set-TextValue 'E25' '  +0.13%  '
Set-TextValue 'E26' '  -3.50%  '
Set-TextValue 'E27' '  -2.67%  '
Set-TextValue 'E28' '  -2.51%  '
Set-TextValue 'D29' '174.80'
Set-TextValue 'E29' '  +3.60%  '
Set-TextValue 'D30' '0.0₃0736'
Set-TextValue 'E30' '  -2.64%  '
Set-TextValue 'E31' '  -2.22%  '
Set-TextValue 'D32' '6.18'
Set-TextValue 'E32' '  -2.81%  '
Set-TextValue 'E33' '  -5.84%  '
Set-TextValue 'E34' '  +0.02%  '
Set-TextValue 'E35' '  -0.24%  '
Set-TextValue 'D36' '17.81'
Set-TextValue 'E36' '  -2.71%  '
Set-TextValue 'E37' '  -6.70%  '
Set-TextValue 'E38' '  -4.66%  '
Set-TextValue 'D39' '36.33'
Set-TextValue 'E39' '  -1.04%  '
Set-TextValue 'E40' '  +2.81%  '
Set-TextValue 'D41' '1.44'
Set-TextValue 'E41' '  -3.86%  '
Set-TextValue 'E42' '  -2.91%  '
Set-TextValue 'D43' '0.583'
Set-TextValue 'E43' '  -3.49%  '
Set-TextValue 'D44' '125.50'
Set-TextValue 'E44' '  -4.90%  '
Set-TextValue 'D45' '259.26'
Set-TextValue 'E45' '  -7.51%  '
Set-TextValue 'E46' '  -6.12%  '
Set-TextValue 'D47' '0.0922'
Set-TextValue 'E47' '  +0.07%  '
Set-TextValue 'E48' '  -2.94%  '
Set-TextValue 'E49' '  -2.53%  '
Set-TextValue 'D50' '17.04'
Set-TextValue 'E50' '  -4.61%  '
Set-TextValue 'D51' '16.30'
Set-TextValue 'E51' '  -4.99%  '
